$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that used to sit around the title
#    ("Searching For Planets In Binary Stars").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Drop the "justify" (both) paragraph alignment from the body paragraphs
#    that currently have it. Setting Alignment to left (0) removes the
#    <w:jc w:val="both"/> element entirely.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Format.Alignment -eq 3) {
        $p.Format.Alignment = 0
    }
}

# ---------------------------------------------------------------------------
# 3. Split the run "Our sample consists of 1100 binary stars, which ..."
#    between the two "11" and "00" halves of "1100", re-inserting the
#    "_GoBack" bookmark (collapsed, start==end) at that split point - this is
#    where an editor's cursor was left after typing "1100" in place of the
#    old number.
# ---------------------------------------------------------------------------
$searchText = "Our sample consists of 1100 binary stars, which represents the largest sample of binary stars for planet search. "
$rng = $d.Content
$found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $runStart = $rng.Start
    $splitPos = $runStart + [string]"Our sample consists of 11".Length

    $bookmarkRange = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)

    # Touch the text right at the split point so the run's serialized
    # "xml:space" gets recomputed (the new left-hand run no longer needs
    # xml:space="preserve" since it has no leading/trailing whitespace).
    $leftRun = $d.Range($runStart, $splitPos)
    $leftRun.InsertAfter("X")
    $marker = $d.Range($splitPos, $splitPos + 1)
    $marker.Text = ""
}
